$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.443826
$ws.Range("H2").Value = 64.331478
$ws.Range("I2").Value = 0.6062978927103765
$ws.Range("J2").Value = 0.6062978927103765
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 32.017979
$ws.Range("N2").Value = 96.05393700000002
$ws.Range("O2").Value = 0.2161524839374964
$ws.Range("P2").Value = 0.2161524839374964
$ws.Range("Q2").Value = 686.5879705476541
$ws.Range("R2").Value = 6179.291734928887
$ws.Range("S2").Value = 0.1310527955154176
$ws.Range("T2").Value = 0.1310527955154176
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.443826
$ws.Range("H3").Value = 64.331478
$ws.Range("I3").Value = 0.6062978927103765
$ws.Range("J3").Value = 0.6062978927103765
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.90901333333333
$ws.Range("N3").Value = 83.72704
$ws.Range("O3").Value = 0.1884129712323413
$ws.Range("P3").Value = 0.1884129712323412
$ws.Range("Q3").Value = 598.47602575168
$ws.Range("R3").Value = 5386.28423176512
$ws.Range("S3").Value = 0.1142343874174693
$ws.Range("T3").Value = 0.1142343874174693
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.443826
$ws.Range("H4").Value = 64.331478
$ws.Range("I4").Value = 0.6062978927103765
$ws.Range("J4").Value = 0.6062978927103765
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.29987433333334
$ws.Range("N4").Value = 60.89962300000001
$ws.Range("O4").Value = 0.1370438858982645
$ws.Range("P4").Value = 0.1370438858982645
$ws.Range("Q4").Value = 435.3069730258661
$ws.Range("R4").Value = 3917.762757232795
$ws.Range("S4").Value = 0.08308941922895906
$ws.Range("T4").Value = 0.08308941922895904
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 21.443826
$ws.Range("H5").Value = 64.331478
$ws.Range("I5").Value = 0.6062978927103765
$ws.Range("J5").Value = 0.6062978927103765
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.468903
$ws.Range("N5").Value = 67.40670900000001
$ws.Range("O5").Value = 0.1516869379794604
$ws.Range("P5").Value = 0.1516869379794604
$ws.Range("Q5").Value = 481.819246342878
$ws.Range("R5").Value = 4336.373217085903
$ws.Range("S5").Value = 0.09196747084863641
$ws.Range("T5").Value = 0.09196747084863639
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 21.443826
$ws.Range("H6").Value = 64.331478
$ws.Range("I6").Value = 0.6062978927103765
$ws.Range("J6").Value = 0.6062978927103765
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 17.57385566666667
$ws.Range("N6").Value = 52.721567
$ws.Range("O6").Value = 0.1186406098495176
$ws.Range("P6").Value = 0.1186406098495176
$ws.Range("Q6").Value = 376.850703065114
$ws.Range("R6").Value = 3391.656327586026
$ws.Range("S6").Value = 0.07193155174163644
$ws.Range("T6").Value = 0.07193155174163643
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 21.443826
$ws.Range("H7").Value = 64.331478
$ws.Range("I7").Value = 0.6062978927103765
$ws.Range("J7").Value = 0.6062978927103765
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 27.85718966666667
$ws.Range("N7").Value = 83.571569
$ws.Range("O7").Value = 0.1880631111029199
$ws.Range("P7").Value = 0.1880631111029199
$ws.Range("Q7").Value = 597.364728060998
$ws.Range("R7").Value = 5376.282552548982
$ws.Range("S7").Value = 0.1140222679582578
$ws.Range("T7").Value = 0.1140222679582578
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.641794
$ws.Range("H8").Value = 10.925382
$ws.Range("I8").Value = 0.1029672609675761
$ws.Range("J8").Value = 0.1029672609675761
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 32.017979
$ws.Range("N8").Value = 96.05393700000002
$ws.Range("O8").Value = 0.2161524839374964
$ws.Range("P8").Value = 0.2161524839374964
$ws.Range("Q8").Value = 116.602883814326
$ws.Range("R8").Value = 1049.425954328934
$ws.Range("S8").Value = 0.02225662922238199
$ws.Range("T8").Value = 0.02225662922238199
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.641794
$ws.Range("H9").Value = 10.925382
$ws.Range("I9").Value = 0.1029672609675761
$ws.Range("J9").Value = 0.1029672609675761
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.90901333333333
$ws.Range("N9").Value = 83.72704
$ws.Range("O9").Value = 0.1884129712323413
$ws.Range("P9").Value = 0.1884129712323412
$ws.Range("Q9").Value = 101.6388773032534
$ws.Range("R9").Value = 914.7498957292801
$ws.Range("S9").Value = 0.01940036757855688
$ws.Range("T9").Value = 0.01940036757855688
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.641794
$ws.Range("H10").Value = 10.925382
$ws.Range("I10").Value = 0.1029672609675761
$ws.Range("J10").Value = 0.1029672609675761
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.29987433333334
$ws.Range("N10").Value = 60.89962300000001
$ws.Range("O10").Value = 0.1370438858982645
$ws.Range("P10").Value = 0.1370438858982645
$ws.Range("Q10").Value = 73.92796054788735
$ws.Range("R10").Value = 665.3516449309861
$ws.Range("S10").Value = 0.01411103356329732
$ws.Range("T10").Value = 0.01411103356329732
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.641794
$ws.Range("H11").Value = 10.925382
$ws.Range("I11").Value = 0.1029672609675761
$ws.Range("J11").Value = 0.1029672609675761
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 22.468903
$ws.Range("N11").Value = 67.40670900000001
$ws.Range("O11").Value = 0.1516869379794604
$ws.Range("P11").Value = 0.1516869379794604
$ws.Range("Q11").Value = 81.82711613198201
$ws.Range("R11").Value = 736.4440451878381
$ws.Range("S11").Value = 0.01561878852830362
$ws.Range("T11").Value = 0.01561878852830362
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.641794
$ws.Range("H12").Value = 10.925382
$ws.Range("I12").Value = 0.1029672609675761
$ws.Range("J12").Value = 0.1029672609675761
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 17.57385566666667
$ws.Range("N12").Value = 52.721567
$ws.Range("O12").Value = 0.1186406098495176
$ws.Range("P12").Value = 0.1186406098495176
$ws.Range("Q12").Value = 64.00036212373267
$ws.Range("R12").Value = 576.003259113594
$ws.Range("S12").Value = 0.01221609863572765
$ws.Range("T12").Value = 0.01221609863572765
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.641794
$ws.Range("H13").Value = 10.925382
$ws.Range("I13").Value = 0.1029672609675761
$ws.Range("J13").Value = 0.1029672609675761
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 27.85718966666667
$ws.Range("N13").Value = 83.571569
$ws.Range("O13").Value = 0.1880631111029199
$ws.Range("P13").Value = 0.1880631111029199
$ws.Range("Q13").Value = 101.4501461849287
$ws.Range("R13").Value = 913.0513156643581
$ws.Range("S13").Value = 0.01936434343930861
$ws.Range("T13").Value = 0.01936434343930861
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 10.28284533333333
$ws.Range("H14").Value = 30.848536
$ws.Range("I14").Value = 0.2907348463220475
$ws.Range("J14").Value = 0.2907348463220475
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 32.017979
$ws.Range("N14").Value = 96.05393700000002
$ws.Range("O14").Value = 0.2161524839374964
$ws.Range("P14").Value = 0.2161524839374964
$ws.Range("Q14").Value = 329.2359259429147
$ws.Range("R14").Value = 2963.123333486233
$ws.Range("S14").Value = 0.06284305919969688
$ws.Range("T14").Value = 0.06284305919969688
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 10.28284533333333
$ws.Range("H15").Value = 30.848536
$ws.Range("I15").Value = 0.2907348463220475
$ws.Range("J15").Value = 0.2907348463220475
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 27.90901333333333
$ws.Range("N15").Value = 83.72704
$ws.Range("O15").Value = 0.1884129712323413
$ws.Range("P15").Value = 0.1884129712323412
$ws.Range("Q15").Value = 286.9840675126044
$ws.Range("R15").Value = 2582.85660761344
$ws.Range("S15").Value = 0.0547782162363151
$ws.Range("T15").Value = 0.0547782162363151
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 10.28284533333333
$ws.Range("H16").Value = 30.848536
$ws.Range("I16").Value = 0.2907348463220475
$ws.Range("J16").Value = 0.2907348463220475
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 20.29987433333334
$ws.Range("N16").Value = 60.89962300000001
$ws.Range("O16").Value = 0.1370438858982645
$ws.Range("P16").Value = 0.1370438858982645
$ws.Range("Q16").Value = 208.7404680557698
$ws.Range("R16").Value = 1878.664212501928
$ws.Range("S16").Value = 0.03984343310600816
$ws.Range("T16").Value = 0.03984343310600814
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 10.28284533333333
$ws.Range("H17").Value = 30.848536
$ws.Range("I17").Value = 0.2907348463220475
$ws.Range("J17").Value = 0.2907348463220475
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 22.468903
$ws.Range("N17").Value = 67.40670900000001
$ws.Range("O17").Value = 0.1516869379794604
$ws.Range("P17").Value = 0.1516869379794604
$ws.Range("Q17").Value = 231.0442543586693
$ws.Range("R17").Value = 2079.398289228024
$ws.Range("S17").Value = 0.04410067860252038
$ws.Range("T17").Value = 0.04410067860252036
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 10.28284533333333
$ws.Range("H18").Value = 30.848536
$ws.Range("I18").Value = 0.2907348463220475
$ws.Range("J18").Value = 0.2907348463220475
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 17.57385566666667
$ws.Range("N18").Value = 52.721567
$ws.Range("O18").Value = 0.1186406098495176
$ws.Range("P18").Value = 0.1186406098495176
$ws.Range("Q18").Value = 180.7092397306569
$ws.Range("R18").Value = 1626.383157575912
$ws.Range("S18").Value = 0.03449295947215349
$ws.Range("T18").Value = 0.03449295947215349
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 10.28284533333333
$ws.Range("H19").Value = 30.848536
$ws.Range("I19").Value = 0.2907348463220475
$ws.Range("J19").Value = 0.2907348463220475
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 27.85718966666667
$ws.Range("N19").Value = 83.571569
$ws.Range("O19").Value = 0.1880631111029199
$ws.Range("P19").Value = 0.1880631111029199
$ws.Range("Q19").Value = 286.4511727636649
$ws.Range("R19").Value = 2578.060554872984
$ws.Range("S19").Value = 0.05467649970535358
$ws.Range("T19").Value = 0.05467649970535358
